{"js": "// Replace the three-digit-by-one-digit multiplication equations in the\n// document's table with the new set of equations, as described by the\n// diff: each old \"A\u00d7B=C\" text run is replaced with a new \"A\u00d7B=C\" text run.\nconst replacements = [\n  [\"494\u00d72=988\", \"277\u00d77=1939\"],\n  [\"800\u00d79=7200\", \"855\u00d78=6840\"],\n  [\"267\u00d73=801\", \"812\u00d74=3248\"],\n  [\"186\u00d78=1488\", \"347\u00d77=2429\"],\n  [\"664\u00d74=2656\", \"964\u00d75=4820\"],\n  [\"721\u00d79=6489\", \"932\u00d78=7456\"],\n  [\"770\u00d72=1540\", \"144\u00d78=1152\"],\n  [\"199\u00d74=796\", \"716\u00d75=3580\"],\n  [\"554\u00d78=4432\", \"235\u00d78=1880\"],\n  [\"542\u00d72=1084\", \"564\u00d77=3948\"],\n  [\"616\u00d74=2464\", \"355\u00d78=2840\"],\n  [\"442\u00d74=1768\", \"449\u00d79=4041\"],\n  [\"590\u00d74=2360\", \"462\u00d79=4158\"],\n  [\"523\u00d77=3661\", \"457\u00d73=1371\"],\n  [\"960\u00d72=1920\", \"256\u00d79=2304\"],\n  [\"393\u00d72=786\", \"342\u00d78=2736\"],\n  [\"693\u00d78=5544\", \"963\u00d75=4815\"],\n  [\"714\u00d75=3570\", \"389\u00d79=3501\"],\n  [\"582\u00d76=3492\", \"691\u00d73=2073\"],\n  [\"156\u00d78=1248\", \"609\u00d78=4872\"],\n  [\"950\u00d75=4750\", \"334\u00d78=2672\"],\n  [\"950\u00d74=3800\", \"702\u00d73=2106\"],\n  [\"190\u00d78=1520\", \"139\u00d76=834\"],\n  [\"849\u00d77=5943\", \"445\u00d77=3115\"],\n  [\"147\u00d78=1176\", \"958\u00d73=2874\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the three-digit-by-one-digit multiplication equations in the\n# document's table with the new set of equations, as described by the\n# diff: each old \"A\u00d7B=C\" text is replaced with a new \"A\u00d7B=C\" text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"494\u00d72=988\", \"277\u00d77=1939\"),\n    @(\"800\u00d79=7200\", \"855\u00d78=6840\"),\n    @(\"267\u00d73=801\", \"812\u00d74=3248\"),\n    @(\"186\u00d78=1488\", \"347\u00d77=2429\"),\n    @(\"664\u00d74=2656\", \"964\u00d75=4820\"),\n    @(\"721\u00d79=6489\", \"932\u00d78=7456\"),\n    @(\"770\u00d72=1540\", \"144\u00d78=1152\"),\n    @(\"199\u00d74=796\", \"716\u00d75=3580\"),\n    @(\"554\u00d78=4432\", \"235\u00d78=1880\"),\n    @(\"542\u00d72=1084\", \"564\u00d77=3948\"),\n    @(\"616\u00d74=2464\", \"355\u00d78=2840\"),\n    @(\"442\u00d74=1768\", \"449\u00d79=4041\"),\n    @(\"590\u00d74=2360\", \"462\u00d79=4158\"),\n    @(\"523\u00d77=3661\", \"457\u00d73=1371\"),\n    @(\"960\u00d72=1920\", \"256\u00d79=2304\"),\n    @(\"393\u00d72=786\", \"342\u00d78=2736\"),\n    @(\"693\u00d78=5544\", \"963\u00d75=4815\"),\n    @(\"714\u00d75=3570\", \"389\u00d79=3501\"),\n    @(\"582\u00d76=3492\", \"691\u00d73=2073\"),\n    @(\"156\u00d78=1248\", \"609\u00d78=4872\"),\n    @(\"950\u00d75=4750\", \"334\u00d78=2672\"),\n    @(\"950\u00d74=3800\", \"702\u00d73=2106\"),\n    @(\"190\u00d78=1520\", \"139\u00d76=834\"),\n    @(\"849\u00d77=5943\", \"445\u00d77=3115\"),\n    @(\"147\u00d78=1176\", \"958\u00d73=2874\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
